$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("A4").Value = 111782529
$ws.Range("B4").Value = 88869
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 2008
$ws.Range("F4").Value = "Fyrflikig jordstjärna"
$ws.Range("G4").Value = "Geastrum quadrifidum"
$ws.Range("H4").Value = "Pers.:Pers."
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "3"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("Q4").Value = 572411
$ws.Range("R4").Value = 6300350
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()

# --- Row 5 ---
$ws.Range("A5").Value = 111782537
$ws.Range("B5").Value = 93388
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2180
$ws.Range("F5").Value = "Blåmossa"
$ws.Range("G5").Value = "Leucobryum glaucum"
$ws.Range("H5").Value = "(Hedw.) Ångstr."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "3"
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("Q5").Value = 572144
$ws.Range("R5").Value = 6300253
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2022-09-12"
$ws.Range("Z5").ClearContents()
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2022-09-12"
$ws.Range("AB5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = 111782531
$ws.Range("B6").Value = 93388
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 2180
$ws.Range("F6").Value = "Blåmossa"
$ws.Range("G6").Value = "Leucobryum glaucum"
$ws.Range("H6").Value = "(Hedw.) Ångstr."
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "10"
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("Q6").Value = 572405
$ws.Range("R6").Value = 6300360
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2022-09-13"
$ws.Range("Z6").ClearContents()
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2022-09-13"
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").Value = "Riklig"
